# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    3  = 1115
    5  = 93
    8  = 11416
    9  = 4323
    13 = 2526
    15 = 123
    16 = 26
    17 = 175
    18 = 500
    19 = 11282
    20 = 11162
    21 = 15
    23 = 12
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    3  = 1115
    5  = 93
    8  = 11416
    9  = 4323
    13 = 2526
    16 = 123
    17 = 26
    18 = 175
    19 = 500
    20 = 11282
    21 = 11162
    22 = 15
    24 = 12
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
